$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.ClearFormatting()
$failures = 0
$ok = $rng.Find.Execute("2023-07-16 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-17 Monday", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 0: 2023-07-16 Sunday -> 2023-07-17 Monday" }
$ok = $rng.Find.Execute("53×83=4399", $true, $false, $false, $false, $false, $true, 1, $false, "75×89=6675", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 1: 53×83=4399 -> 75×89=6675" }
$ok = $rng.Find.Execute("76×15=1140", $true, $false, $false, $false, $false, $true, 1, $false, "14×25=350", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 2: 76×15=1140 -> 14×25=350" }
$ok = $rng.Find.Execute("10×93=930", $true, $false, $false, $false, $false, $true, 1, $false, "85×73=6205", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 3: 10×93=930 -> 85×73=6205" }
$ok = $rng.Find.Execute("93×100=9300", $true, $false, $false, $false, $false, $true, 1, $false, "42×63=2646", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 4: 93×100=9300 -> 42×63=2646" }
$ok = $rng.Find.Execute("53×70=3710", $true, $false, $false, $false, $false, $true, 1, $false, "48×41=1968", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 5: 53×70=3710 -> 48×41=1968" }
$ok = $rng.Find.Execute("50×29=1450", $true, $false, $false, $false, $false, $true, 1, $false, "54×49=2646", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 6: 50×29=1450 -> 54×49=2646" }
$ok = $rng.Find.Execute("73×43=3139", $true, $false, $false, $false, $false, $true, 1, $false, "60×56=3360", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 7: 73×43=3139 -> 60×56=3360" }
$ok = $rng.Find.Execute("53×30=1590", $true, $false, $false, $false, $false, $true, 1, $false, "67×80=5360", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 8: 53×30=1590 -> 67×80=5360" }
$ok = $rng.Find.Execute("76×10=760", $true, $false, $false, $false, $false, $true, 1, $false, "83×35=2905", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 9: 76×10=760 -> 83×35=2905" }
$ok = $rng.Find.Execute("45×31=1395", $true, $false, $false, $false, $false, $true, 1, $false, "41×94=3854", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 10: 45×31=1395 -> 41×94=3854" }
$ok = $rng.Find.Execute("73×43=3139", $true, $false, $false, $false, $false, $true, 1, $false, "58×80=4640", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 11: 73×43=3139 -> 58×80=4640" }
$ok = $rng.Find.Execute("34×34=1156", $true, $false, $false, $false, $false, $true, 1, $false, "15×84=1260", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 12: 34×34=1156 -> 15×84=1260" }
$ok = $rng.Find.Execute("29×37=1073", $true, $false, $false, $false, $false, $true, 1, $false, "27×59=1593", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 13: 29×37=1073 -> 27×59=1593" }
$ok = $rng.Find.Execute("42×21=882", $true, $false, $false, $false, $false, $true, 1, $false, "28×72=2016", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 14: 42×21=882 -> 28×72=2016" }
$ok = $rng.Find.Execute("29×86=2494", $true, $false, $false, $false, $false, $true, 1, $false, "37×84=3108", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 15: 29×86=2494 -> 37×84=3108" }
$ok = $rng.Find.Execute("72×29=2088", $true, $false, $false, $false, $false, $true, 1, $false, "35×41=1435", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 16: 72×29=2088 -> 35×41=1435" }
$ok = $rng.Find.Execute("60×94=5640", $true, $false, $false, $false, $false, $true, 1, $false, "33×44=1452", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 17: 60×94=5640 -> 33×44=1452" }
$ok = $rng.Find.Execute("82×68=5576", $true, $false, $false, $false, $false, $true, 1, $false, "43×36=1548", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 18: 82×68=5576 -> 43×36=1548" }
$ok = $rng.Find.Execute("37×23=851", $true, $false, $false, $false, $false, $true, 1, $false, "42×63=2646", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 19: 37×23=851 -> 42×63=2646" }
$ok = $rng.Find.Execute("70×39=2730", $true, $false, $false, $false, $false, $true, 1, $false, "13×51=663", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 20: 70×39=2730 -> 13×51=663" }
$ok = $rng.Find.Execute("45×25=1125", $true, $false, $false, $false, $false, $true, 1, $false, "75×93=6975", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 21: 45×25=1125 -> 75×93=6975" }
$ok = $rng.Find.Execute("11×60=660", $true, $false, $false, $false, $false, $true, 1, $false, "35×50=1750", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 22: 11×60=660 -> 35×50=1750" }
$ok = $rng.Find.Execute("18×36=648", $true, $false, $false, $false, $false, $true, 1, $false, "46×94=4324", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 23: 18×36=648 -> 46×94=4324" }
$ok = $rng.Find.Execute("98×82=8036", $true, $false, $false, $false, $false, $true, 1, $false, "31×40=1240", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 24: 98×82=8036 -> 31×40=1240" }
$ok = $rng.Find.Execute("41×70=2870", $true, $false, $false, $false, $false, $true, 1, $false, "40×21=840", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 25: 41×70=2870 -> 40×21=840" }
$ok = $rng.Find.Execute("83×100=8300", $true, $false, $false, $false, $false, $true, 1, $false, "18×94=1692", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 26: 83×100=8300 -> 18×94=1692" }
$ok = $rng.Find.Execute("86×29=2494", $true, $false, $false, $false, $false, $true, 1, $false, "82×58=4756", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 27: 86×29=2494 -> 82×58=4756" }
$ok = $rng.Find.Execute("97×52=5044", $true, $false, $false, $false, $false, $true, 1, $false, "77×65=5005", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 28: 97×52=5044 -> 77×65=5005" }
$ok = $rng.Find.Execute("38×39=1482", $true, $false, $false, $false, $false, $true, 1, $false, "30×30=900", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 29: 38×39=1482 -> 30×30=900" }
$ok = $rng.Find.Execute("72×27=1944", $true, $false, $false, $false, $false, $true, 1, $false, "78×12=936", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 30: 72×27=1944 -> 78×12=936" }
$ok = $rng.Find.Execute("61×30=1830", $true, $false, $false, $false, $false, $true, 1, $false, "37×64=2368", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 31: 61×30=1830 -> 37×64=2368" }
$ok = $rng.Find.Execute("74×85=6290", $true, $false, $false, $false, $false, $true, 1, $false, "11×64=704", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 32: 74×85=6290 -> 11×64=704" }
$ok = $rng.Find.Execute("90×11=990", $true, $false, $false, $false, $false, $true, 1, $false, "84×29=2436", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 33: 90×11=990 -> 84×29=2436" }
$ok = $rng.Find.Execute("59×27=1593", $true, $false, $false, $false, $false, $true, 1, $false, "37×96=3552", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 34: 59×27=1593 -> 37×96=3552" }
$ok = $rng.Find.Execute("11×67=737", $true, $false, $false, $false, $false, $true, 1, $false, "27×19=513", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 35: 11×67=737 -> 27×19=513" }
$ok = $rng.Find.Execute("75×100=7500", $true, $false, $false, $false, $false, $true, 1, $false, "57×53=3021", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 36: 75×100=7500 -> 57×53=3021" }
$ok = $rng.Find.Execute("84×93=7812", $true, $false, $false, $false, $false, $true, 1, $false, "47×68=3196", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 37: 84×93=7812 -> 47×68=3196" }
$ok = $rng.Find.Execute("59×12=708", $true, $false, $false, $false, $false, $true, 1, $false, "20×97=1940", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 38: 59×12=708 -> 20×97=1940" }
$ok = $rng.Find.Execute("54×19=1026", $true, $false, $false, $false, $false, $true, 1, $false, "72×28=2016", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 39: 54×19=1026 -> 72×28=2016" }
$ok = $rng.Find.Execute("50×57=2850", $true, $false, $false, $false, $false, $true, 1, $false, "99×35=3465", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 40: 50×57=2850 -> 99×35=3465" }
$ok = $rng.Find.Execute("92×48=4416", $true, $false, $false, $false, $false, $true, 1, $false, "44×18=792", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 41: 92×48=4416 -> 44×18=792" }
$ok = $rng.Find.Execute("90×84=7560", $true, $false, $false, $false, $false, $true, 1, $false, "62×58=3596", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 42: 90×84=7560 -> 62×58=3596" }
$ok = $rng.Find.Execute("22×49=1078", $true, $false, $false, $false, $false, $true, 1, $false, "58×76=4408", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 43: 22×49=1078 -> 58×76=4408" }
$ok = $rng.Find.Execute("93×25=2325", $true, $false, $false, $false, $false, $true, 1, $false, "82×81=6642", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 44: 93×25=2325 -> 82×81=6642" }
$ok = $rng.Find.Execute("40×72=2880", $true, $false, $false, $false, $false, $true, 1, $false, "47×51=2397", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 45: 40×72=2880 -> 47×51=2397" }
$ok = $rng.Find.Execute("62×86=5332", $true, $false, $false, $false, $false, $true, 1, $false, "70×31=2170", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 46: 62×86=5332 -> 70×31=2170" }
$ok = $rng.Find.Execute("65×29=1885", $true, $false, $false, $false, $false, $true, 1, $false, "58×75=4350", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 47: 65×29=1885 -> 58×75=4350" }
$ok = $rng.Find.Execute("71×57=4047", $true, $false, $false, $false, $false, $true, 1, $false, "98×83=8134", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 48: 71×57=4047 -> 98×83=8134" }
$ok = $rng.Find.Execute("56×48=2688", $true, $false, $false, $false, $false, $true, 1, $false, "96×37=3552", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 49: 56×48=2688 -> 96×37=3552" }
$ok = $rng.Find.Execute("73×93=6789", $true, $false, $false, $false, $false, $true, 1, $false, "18×93=1674", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 50: 73×93=6789 -> 18×93=1674" }
$ok = $rng.Find.Execute("10×97=970", $true, $false, $false, $false, $false, $true, 1, $false, "24×67=1608", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 51: 10×97=970 -> 24×67=1608" }
$ok = $rng.Find.Execute("85×20=1700", $true, $false, $false, $false, $false, $true, 1, $false, "58×40=2320", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 52: 85×20=1700 -> 58×40=2320" }
$ok = $rng.Find.Execute("58×57=3306", $true, $false, $false, $false, $false, $true, 1, $false, "88×86=7568", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 53: 58×57=3306 -> 88×86=7568" }
$ok = $rng.Find.Execute("73×77=5621", $true, $false, $false, $false, $false, $true, 1, $false, "11×85=935", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 54: 73×77=5621 -> 11×85=935" }
$ok = $rng.Find.Execute("22×82=1804", $true, $false, $false, $false, $false, $true, 1, $false, "53×93=4929", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 55: 22×82=1804 -> 53×93=4929" }
$ok = $rng.Find.Execute("17×25=425", $true, $false, $false, $false, $false, $true, 1, $false, "15×15=225", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 56: 17×25=425 -> 15×15=225" }
$ok = $rng.Find.Execute("75×34=2550", $true, $false, $false, $false, $false, $true, 1, $false, "28×75=2100", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 57: 75×34=2550 -> 28×75=2100" }
$ok = $rng.Find.Execute("22×14=308", $true, $false, $false, $false, $false, $true, 1, $false, "27×15=405", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 58: 22×14=308 -> 27×15=405" }
$ok = $rng.Find.Execute("26×75=1950", $true, $false, $false, $false, $false, $true, 1, $false, "35×17=595", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 59: 26×75=1950 -> 35×17=595" }
$ok = $rng.Find.Execute("83×47=3901", $true, $false, $false, $false, $false, $true, 1, $false, "15×66=990", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 60: 83×47=3901 -> 15×66=990" }
$ok = $rng.Find.Execute("60×93=5580", $true, $false, $false, $false, $false, $true, 1, $false, "56×30=1680", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 61: 60×93=5580 -> 56×30=1680" }
$ok = $rng.Find.Execute("77×24=1848", $true, $false, $false, $false, $false, $true, 1, $false, "99×93=9207", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 62: 77×24=1848 -> 99×93=9207" }
$ok = $rng.Find.Execute("44×58=2552", $true, $false, $false, $false, $false, $true, 1, $false, "18×19=342", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 63: 44×58=2552 -> 18×19=342" }
$ok = $rng.Find.Execute("68×29=1972", $true, $false, $false, $false, $false, $true, 1, $false, "35×57=1995", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 64: 68×29=1972 -> 35×57=1995" }
$ok = $rng.Find.Execute("62×22=1364", $true, $false, $false, $false, $false, $true, 1, $false, "63×42=2646", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 65: 62×22=1364 -> 63×42=2646" }
$ok = $rng.Find.Execute("91×80=7280", $true, $false, $false, $false, $false, $true, 1, $false, "48×85=4080", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 66: 91×80=7280 -> 48×85=4080" }
$ok = $rng.Find.Execute("12×39=468", $true, $false, $false, $false, $false, $true, 1, $false, "97×16=1552", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 67: 12×39=468 -> 97×16=1552" }
$ok = $rng.Find.Execute("86×71=6106", $true, $false, $false, $false, $false, $true, 1, $false, "53×24=1272", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 68: 86×71=6106 -> 53×24=1272" }
$ok = $rng.Find.Execute("45×17=765", $true, $false, $false, $false, $false, $true, 1, $false, "64×43=2752", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 69: 45×17=765 -> 64×43=2752" }
$ok = $rng.Find.Execute("16×69=1104", $true, $false, $false, $false, $false, $true, 1, $false, "32×39=1248", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 70: 16×69=1104 -> 32×39=1248" }
$ok = $rng.Find.Execute("41×52=2132", $true, $false, $false, $false, $false, $true, 1, $false, "80×53=4240", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 71: 41×52=2132 -> 80×53=4240" }
$ok = $rng.Find.Execute("85×35=2975", $true, $false, $false, $false, $false, $true, 1, $false, "72×43=3096", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 72: 85×35=2975 -> 72×43=3096" }
$ok = $rng.Find.Execute("28×85=2380", $true, $false, $false, $false, $false, $true, 1, $false, "95×95=9025", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 73: 28×85=2380 -> 95×95=9025" }
$ok = $rng.Find.Execute("23×97=2231", $true, $false, $false, $false, $false, $true, 1, $false, "86×57=4902", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 74: 23×97=2231 -> 86×57=4902" }
$ok = $rng.Find.Execute("81×31=2511", $true, $false, $false, $false, $false, $true, 1, $false, "58×22=1276", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 75: 81×31=2511 -> 58×22=1276" }
$ok = $rng.Find.Execute("11×14=154", $true, $false, $false, $false, $false, $true, 1, $false, "55×26=1430", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 76: 11×14=154 -> 55×26=1430" }
$ok = $rng.Find.Execute("63×45=2835", $true, $false, $false, $false, $false, $true, 1, $false, "50×63=3150", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 77: 63×45=2835 -> 50×63=3150" }
$ok = $rng.Find.Execute("35×24=840", $true, $false, $false, $false, $false, $true, 1, $false, "42×92=3864", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 78: 35×24=840 -> 42×92=3864" }
$ok = $rng.Find.Execute("53×62=3286", $true, $false, $false, $false, $false, $true, 1, $false, "94×11=1034", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 79: 53×62=3286 -> 94×11=1034" }
$ok = $rng.Find.Execute("97×81=7857", $true, $false, $false, $false, $false, $true, 1, $false, "37×11=407", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 80: 97×81=7857 -> 37×11=407" }
$ok = $rng.Find.Execute("59×57=3363", $true, $false, $false, $false, $false, $true, 1, $false, "23×20=460", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 81: 59×57=3363 -> 23×20=460" }
$ok = $rng.Find.Execute("96×41=3936", $true, $false, $false, $false, $false, $true, 1, $false, "50×45=2250", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 82: 96×41=3936 -> 50×45=2250" }
$ok = $rng.Find.Execute("97×53=5141", $true, $false, $false, $false, $false, $true, 1, $false, "14×46=644", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 83: 97×53=5141 -> 14×46=644" }
$ok = $rng.Find.Execute("76×90=6840", $true, $false, $false, $false, $false, $true, 1, $false, "43×42=1806", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 84: 76×90=6840 -> 43×42=1806" }
$ok = $rng.Find.Execute("77×49=3773", $true, $false, $false, $false, $false, $true, 1, $false, "64×16=1024", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 85: 77×49=3773 -> 64×16=1024" }
$ok = $rng.Find.Execute("62×16=992", $true, $false, $false, $false, $false, $true, 1, $false, "59×18=1062", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 86: 62×16=992 -> 59×18=1062" }
$ok = $rng.Find.Execute("13×27=351", $true, $false, $false, $false, $false, $true, 1, $false, "34×63=2142", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 87: 13×27=351 -> 34×63=2142" }
$ok = $rng.Find.Execute("89×23=2047", $true, $false, $false, $false, $false, $true, 1, $false, "70×89=6230", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 88: 89×23=2047 -> 70×89=6230" }
$ok = $rng.Find.Execute("42×82=3444", $true, $false, $false, $false, $false, $true, 1, $false, "42×16=672", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 89: 42×82=3444 -> 42×16=672" }
$ok = $rng.Find.Execute("73×25=1825", $true, $false, $false, $false, $false, $true, 1, $false, "54×84=4536", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 90: 73×25=1825 -> 54×84=4536" }
$ok = $rng.Find.Execute("100×100=10000", $true, $false, $false, $false, $false, $true, 1, $false, "36×92=3312", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 91: 100×100=10000 -> 36×92=3312" }
$ok = $rng.Find.Execute("33×42=1386", $true, $false, $false, $false, $false, $true, 1, $false, "69×17=1173", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 92: 33×42=1386 -> 69×17=1173" }
$ok = $rng.Find.Execute("36×96=3456", $true, $false, $false, $false, $false, $true, 1, $false, "14×11=154", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 93: 36×96=3456 -> 14×11=154" }
$ok = $rng.Find.Execute("50×83=4150", $true, $false, $false, $false, $false, $true, 1, $false, "86×37=3182", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 94: 50×83=4150 -> 86×37=3182" }
$ok = $rng.Find.Execute("63×33=2079", $true, $false, $false, $false, $false, $true, 1, $false, "15×84=1260", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 95: 63×33=2079 -> 15×84=1260" }
$ok = $rng.Find.Execute("10×47=470", $true, $false, $false, $false, $false, $true, 1, $false, "28×12=336", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 96: 10×47=470 -> 28×12=336" }
$ok = $rng.Find.Execute("54×87=4698", $true, $false, $false, $false, $false, $true, 1, $false, "46×45=2070", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 97: 54×87=4698 -> 46×45=2070" }
$ok = $rng.Find.Execute("68×55=3740", $true, $false, $false, $false, $false, $true, 1, $false, "54×85=4590", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 98: 68×55=3740 -> 54×85=4590" }
$ok = $rng.Find.Execute("94×33=3102", $true, $false, $false, $false, $false, $true, 1, $false, "52×74=3848", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 99: 94×33=3102 -> 52×74=3848" }
$ok = $rng.Find.Execute("60×72=4320", $true, $false, $false, $false, $false, $true, 1, $false, "14×43=602", 1)
if (-not $ok) { $failures++; Write-Output "FAILED at index 100: 60×72=4320 -> 14×43=602" }
Write-Output "Done. failures=$failures"
